$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Add the new "People" worksheet right after "TypeAhesd" and make it the
#    active sheet (matches activeTab="1" / tabSelected moving to sheet2).
# ---------------------------------------------------------------------------
$people = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet1)
$people.Name = "People"

# Header row (row 1) - reuses the same header labels/styles as "TypeAhesd".
$people.Range("A1").Value = "TESTNAME"
$people.Range("B1").Value = "DESCRIPTION"
$people.Range("C1").Value = "HOST"
$people.Range("D1").Value = "API"
$people.Range("E1").Value = "METHOD"
$people.Range("F1").Value = "HEADERS"
$people.Range("G1").Value = "QUERYSTRING"
$people.Range("H1").Value = "BODY"
$people.Range("I1").Value = "DEPENDENCYTESTS"
$people.Range("J1").Value = "VALIDATIONS"
$people.Range("K1").Value = "STORE"
$people.Range("L1").Value = "STATUS"

$people.Range("A1:L1").Font.Bold = $true
$people.Range("A1:L1").Interior.Color = 65535
$people.Range("A1:L1").Borders.LineStyle = 1
$people.Range("A1:L1").HorizontalAlignment = -4108
$people.Range("H1").WrapText = $true
$people.Range("J1").WrapText = $true

# Data rows - fill column by column (A, B, G) for every row first, then
# go back and fill in the J (VALIDATIONS) column, mirroring the order the
# strings were authored in (keeps shared-string ordering identical).
$people.Range("A2").Value = "OPQA_1222"
$people.Range("B2").Value = "Verify that Type Ahead returns peoples by passing  query."
$people.Range("G2").Value = "?query=projec&source=people"

$people.Range("A3").Value = "OPQA_1222"
$people.Range("B3").Value = "Verify that Type Ahead returns peoples by passing  user first name."
$people.Range("G3").Value = "?query=project&source=people"

$people.Range("A4").Value = "OPQA_1222"
$people.Range("B4").Value = "Verify that Type Ahead returns peoples by passing  user last name."
$people.Range("G4").Value = "?query=Neon1&source=people"

$people.Range("J2").Value = "status=200||source=people||suggestions.keyword=projec||suggestions.info.value=Project Neon1||suggestions.info.value=Project Neon2||suggestions.info.value=Project Neon3"
$people.Range("J3").Value = "status=200||source=people||suggestions.keyword=project||suggestions.info.value=Project Neon1||suggestions.info.value=Project Neon2||suggestions.info.value=Project Neon3"
$people.Range("J4").Value = "status=200||source=people||suggestions.keyword=Neon1||suggestions.info.value=Project Neon1"

# Remaining columns (reuse existing shared strings - "1PTYPEAHEAD", "/suggest", "GET")
$people.Range("C2").Value = "1PTYPEAHEAD"
$people.Range("D2").Value = "/suggest"
$people.Range("E2").Value = "GET"
$people.Range("H2").Style = "Normal"

$people.Range("C3").Value = "1PTYPEAHEAD"
$people.Range("D3").Value = "/suggest"
$people.Range("E3").Value = "GET"
$people.Range("H3").Style = "Normal"

$people.Range("C4").Value = "1PTYPEAHEAD"
$people.Range("D4").Value = "/suggest"
$people.Range("E4").Value = "GET"
$people.Range("H4").Style = "Normal"

# Formatting to match "TypeAhesd": wrap text on description/validation cells.
$people.Range("B2:B4").WrapText = $true
$people.Range("J2:J4").WrapText = $true

$people.Rows.Item(2).RowHeight = 60
$people.Rows.Item(3).RowHeight = 60
$people.Rows.Item(4).RowHeight = 30

# Column widths (approximate the original authored best-fit widths).
$people.Columns.Item(1).ColumnWidth = 10.5
$people.Columns.Item(2).ColumnWidth = 38.5
$people.Columns.Item(3).ColumnWidth = 12.75
$people.Columns.Item(4).ColumnWidth = 43
$people.Columns.Item(5).ColumnWidth = 8
$people.Columns.Item(6).ColumnWidth = 21
$people.Columns.Item(7).ColumnWidth = 30
$people.Columns.Item(8).ColumnWidth = 17
$people.Columns.Item(9).ColumnWidth = 18
$people.Columns.Item(10).ColumnWidth = 50
$people.Columns.Item(11).ColumnWidth = 23
$people.Columns.Item(12).ColumnWidth = 10.5

$people.Range("E1").Select()
$people.Range("L2:L4").Select()

Write-Output "done"
